$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 8: set B8 to new value, H8/I8 to attributeType/AttributeType
$ws.Range("B8").Value = "AttributeTypeAndFloat"
$ws.Range("H8").Value = "attributeType"
$ws.Range("I8").Value = "AttributeType"

# Row 9: H9/I9 value/float
$ws.Range("H9").Value = "value"
$ws.Range("I9").Value = "float"

# Selection change
$ws.Range("F17").Select()

# Column B widened to fit the longer "AttributeTypeAndFloat" label
# (splits the old B:C shared width so column B gets its own best-fit width)
$ws.Columns.Item(2).ColumnWidth = 22.140625
